$d = $word.ActiveDocument

# --- Paragraph 1: "Autor" / "es " / ":" / " " (bold, 4 runs) + nbsp run ---
# Locate the paragraph that starts the "Autores : " label (first occurrence,
# which currently has its bold text split across separate runs: "Autor",
# "es ", ":", " "). Merge those runs into a single "Autores : " run via
# Find/Replace (same text in, forces a run-consolidation), then append the
# new "Mario y Diego" text (unbolded) right after the existing trailing
# non-breaking-space run.
$nbsp = [char]0x00A0
$target1 = "Autores : " + $nbsp + "`r"
$target2 = "Autor:" + $nbsp + "`r"

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text

    if ($t -eq $target1) {
        $r = $p.Range
        $null = $r.Find.Execute("Autores : ", $true, $false, $false, $false, $false, $true, 1, $false, "Autores : ", 2)

        $p2 = $d.Paragraphs.Item($i)
        $r2 = $p2.Range
        $pos = $r2.End - 1
        $ins = $d.Range($pos, $pos)
        $ins.InsertAfter("Mario y Diego")
    }
    elseif ($t -eq $target2) {
        $pos = $p.Range.End - 1
        $ins = $d.Range($pos, $pos)
        $ins.InsertAfter(" Mario y Diego")
    }
}
